# 20250828 Budget file update
# Bump the Oct/Nov/Dec 2025 (columns BS/BT/BU) budget figures on rows 2-4
# from 20000 to 21000, and leave the selection where the editor left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("BS2:BU2").Value = 21000
$ws.Range("BS3:BU3").Value = 21000
$ws.Range("BS4:BU4").Value = 21000

$ws.Activate()
$ws.Range("BS21").Select()
